$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows: row#, A (Sending cluster), B (Ligand symbol), C (Receptor symbol), D (Target cluster),
# E..T numeric metrics
$rows = @(
    @(2, "ECs", "Wnt9a", "Fzd4", "ECs", 2, 0.6666666666666666, 0.2740863333333333, 0.822259, 0.05829365131129327, 0.05829365131129326, 3, 1, 31.38723566666667, 94.161707, 0.5539598599114094, 0.5539598599114095, 8.602812337345888, 77.425311036113, 0.03229234291412857, 0.03229234291412857),
    @(3, "ECs", "Wnt9a", "Fzd4", "FAPs", 2, 0.6666666666666666, 0.2740863333333333, 0.822259, 0.05829365131129327, 0.05829365131129326, 3, 1, 17.55525033333333, 52.665751, 0.3098362697066353, 0.3098362697066353, 4.811654194612111, 43.304887751509, 0.01806148746987042, 0.01806148746987041),
    @(4, "ECs", "Wnt9a", "Fzd4", "sCs", 2, 0.6666666666666666, 0.2740863333333333, 0.822259, 0.05829365131129327, 0.05829365131129326, 3, 1, 7.717279333333334, 23.151838, 0.1362038703819552, 0.1362038703819552, 2.115200795782445, 19.036807162042, 0.007939820927294283, 0.007939820927294282),
    @(5, "FAPs", "Wnt9a", "Fzd4", "ECs", 3, 1, 3.335609, 10.006827, 0.7094291261882629, 0.7094291261882628, 3, 1, 31.38723566666667, 94.161707, 0.5539598599114094, 0.5539598599114095, 104.6955457748543, 942.259911973689, 0.3929952593603237, 0.3929952593603237),
    @(6, "FAPs", "Wnt9a", "Fzd4", "FAPs", 3, 1, 3.335609, 10.006827, 0.7094291261882629, 0.7094291261882628, 3, 1, 17.55525033333333, 52.665751, 0.3098362697066353, 0.3098362697066353, 58.55745100911967, 527.017059082077, 0.2198068740794092, 0.2198068740794092),
    @(7, "FAPs", "Wnt9a", "Fzd4", "sCs", 3, 1, 3.335609, 10.006827, 0.7094291261882629, 0.7094291261882628, 3, 1, 7.717279333333334, 23.151838, 0.1362038703819552, 0.1362038703819552, 25.74182639978067, 231.676437598026, 0.09662699274852991, 0.0966269927485299),
    @(8, "sCs", "Wnt9a", "Fzd4", "ECs", 3, 1, 1.092126, 3.276378, 0.2322772225004438, 0.2322772225004438, 3, 1, 31.38723566666667, 94.161707, 0.5539598599114094, 0.5539598599114095, 34.27881613969401, 308.509345257246, 0.1286722576369571, 0.1286722576369571),
    @(9, "sCs", "Wnt9a", "Fzd4", "FAPs", 3, 1, 1.092126, 3.276378, 0.2322772225004438, 0.2322772225004438, 3, 1, 17.55525033333333, 52.665751, 0.3098362697066353, 0.3098362697066353, 19.172545325542, 172.552907929878, 0.07196790815735565, 0.07196790815735564),
    @(10, "sCs", "Wnt9a", "Fzd4", "sCs", 3, 1, 1.092126, 3.276378, 0.2322772225004438, 0.2322772225004438, 3, 1, 7.717279333333334, 23.151838, 0.1362038703819552, 0.1362038703819552, 8.428241409196001, 75.85417268276402, 0.03163705670613103, 0.03163705670613102)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    for ($c = 5; $c -le 20; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c]
    }
}

Write-Host "updated rows 2-10"
